$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the first
#    paragraph (the "Play Free Auspicious Fortune God Slot" Heading1).
#    We source the new paragraph's run structure from the existing
#    bold-formatted "Play Free Auspicious Fortune God Slot" paragraph
#    further down in the document (it already has the leading empty
#    run + a run carrying <w:rPr><w:b/></w:rPr>, which is exactly the
#    shape we need), paste a copy of it in place, then swap in the
#    new text.
# ------------------------------------------------------------------

$origCount = $d.Paragraphs.Count
$boldSourceIndex = $origCount - 1

$srcPara = $d.Paragraphs.Item($boldSourceIndex)
$srcPara.Range.Copy()

$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"
$metaPara.Range.Paste()

$boldStart = $metaPara.Range.Start
$boldEnd = $metaPara.Range.End - 1
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Text = "Meta description"

$afterBold = $d.Range($metaPara.Range.End - 1, $metaPara.Range.End - 1)
$afterBold.InsertAfter(": Read our review of Auspicious Fortune God, a Chinese-themed slot game with Expanding Wilds, Respins, and Customization Options. Play for free.")

# ------------------------------------------------------------------
# 2) Remove the old bold "Play Free Auspicious Fortune God Slot"
#    paragraph that used to sit right before the meta-description
#    (italic) paragraph near the end of the document, and rewrite the
#    italic paragraph's text into the new image-generation prompt.
#    Every paragraph index from step 1 onward has shifted by +1
#    because of the paragraph we inserted above.
# ------------------------------------------------------------------

$tailBoldIndex = $boldSourceIndex + 1
$tailItalicIndex = $tailBoldIndex + 1

$d.Paragraphs.Item($tailBoldIndex).Range.Delete()

# After the delete, the italic paragraph shifts back down by one.
$italicPara = $d.Paragraphs.Item($tailBoldIndex)
$italicStart = $italicPara.Range.Start
$italicEnd = $italicPara.Range.End - 1
$italicRange = $d.Range($italicStart, $italicEnd)
$italicRange.Text = "Create an eye-catching feature image for ""Auspicious Fortune God"" that showcases a happy Maya warrior with glasses in cartoon style. Use bright and bold colors that are reminiscent of the game's Chinese-theme, such as red and gold, to catch the viewer's attention. Position the warrior in a confident and charismatic pose, with his hands on his hips and a big smile on his face. Make sure to incorporate the game's logo into the image and any other relevant symbols such as fortune deities and money trees. The overall image should convey a sense of excitement and fun while highlighting the game's unique features."
